$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row at position 3 (shifts existing rows 3-38 down to 4-39) ---
$ws.Rows(3).Insert()

# Row 3
$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = '6262'
$ws.Range("A3").Style = "Normal"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = '6/25/2025'
$ws.Range("B3").Style = "Normal"
$ws.Range("C3").Value = 'MIGUELETES 1330'
$ws.Range("D3").Value = 14
$ws.Range("E3").Value = 'ICD30465943'
$ws.Range("F3").Value = 'Optical Power'
$ws.Range("G3").Value = 'Pendiente'
$ws.Range("H3").Value = 'Cables en panza'
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = '{"direccionesNormalizadas": [{"altura": 1330, "cod_calle": 13079, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.440291", "y": "-34.562841"}, "direccion": "MIGUELETES 1330, CABA", "nombre_calle": "MIGUELETES", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K3").Value = -58.440291
$ws.Range("L3").Value = -34.562841
$ws.Range("M3").Value = 'Colegiales'
$ws.Range("N3").Value = 'Capital Norte'

# --- Append 15 brand-new rows after the (shifted) last existing row (now row 39) ---

# Row 40
$ws.Range("A40").NumberFormat = "@"
$ws.Range("A40").Value = '7269'
$ws.Range("A40").Style = "Normal"
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = '9/22/2025'
$ws.Range("B40").Style = "Normal"
$ws.Range("C40").Value = 'BACACAY 2936'
$ws.Range("D40").Value = 7
$ws.Range("E40").Value = 'ICD30951858'
$ws.Range("F40").Value = 'Optical Power'
$ws.Range("G40").Value = 'Pendiente'
$ws.Range("H40").Value = 'Tendido a baja altura'
$ws.Range("I40").Value = 1
$ws.Range("J40").Value = '{"direccionesNormalizadas": [{"altura": 2936, "cod_calle": 2002, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.471106", "y": "-34.628593"}, "direccion": "BACACAY 2936, CABA", "nombre_calle": "BACACAY", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K40").Value = -58.471106
$ws.Range("L40").Value = -34.628593
$ws.Range("M40").Value = 'Devoto'
$ws.Range("N40").Value = 'Capital Norte'

# Row 41
$ws.Range("A41").NumberFormat = "@"
$ws.Range("A41").Value = '7276'
$ws.Range("A41").Style = "Normal"
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = '9/22/2025'
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").Value = 'LA FRONDA 1684'
$ws.Range("D41").Value = 11
$ws.Range("E41").Value = 'ICD30952118'
$ws.Range("F41").Value = 'Optical Power'
$ws.Range("G41").Value = 'Pendiente'
$ws.Range("H41").Value = 'Baja de cliente a baja altura y agarrada a la columna'
$ws.Range("I41").Value = 1
$ws.Range("J41").Value = '{"direccionesNormalizadas": [{"altura": 1684, "cod_calle": 12022, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.469037", "y": "-34.610941"}, "direccion": "LA FRONDA 1684, CABA", "nombre_calle": "LA FRONDA", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K41").Value = -58.469037
$ws.Range("L41").Value = -34.610941
$ws.Range("M41").Value = 'Paternal'
$ws.Range("N41").Value = 'Capital Norte'

# Row 42
$ws.Range("A42").NumberFormat = "@"
$ws.Range("A42").Value = '7279'
$ws.Range("A42").Style = "Normal"
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = '9/22/2025'
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = 'BACACAY 915'
$ws.Range("D42").Value = 6
$ws.Range("E42").Value = 'Pendiente ADM'
$ws.Range("F42").Value = 'Optical Power'
$ws.Range("G42").Value = 'Pendiente'
$ws.Range("H42").Value = 'Cable cortados y en panza'
$ws.Range("I42").Value = 1
$ws.Range("J42").Value = '{"direccionesNormalizadas": [{"altura": 915, "cod_calle": 2002, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.442669", "y": "-34.618300"}, "direccion": "BACACAY 915, CABA", "nombre_calle": "BACACAY", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K42").Value = -58.442669
$ws.Range("L42").Value = -34.6183
$ws.Range("M42").Value = 'Almagro'
$ws.Range("N42").Value = 'Capital Sur'

# Row 43
$ws.Range("A43").NumberFormat = "@"
$ws.Range("A43").Value = '7280'
$ws.Range("A43").Style = "Normal"
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = '9/22/2025'
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = 'ORTEGA 974'
$ws.Range("D43").Value = 6
$ws.Range("E43").Value = 'Pendiente ADM'
$ws.Range("F43").Value = 'Optical Power'
$ws.Range("G43").Value = 'Pendiente'
$ws.Range("H43").Value = 'Cable en panza y cortado'
$ws.Range("I43").Value = 1
$ws.Range("J43").Value = '{"direccionesNormalizadas": [{"altura": 974, "cod_calle": 16031, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.445007", "y": "-34.614987"}, "direccion": "ORTEGA 974, CABA", "nombre_calle": "ORTEGA", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K43").Value = -58.445007
$ws.Range("L43").Value = -34.614987
$ws.Range("M43").Value = 'Almagro'
$ws.Range("N43").Value = 'Capital Sur'

# Row 44
$ws.Range("A44").NumberFormat = "@"
$ws.Range("A44").Value = '7281'
$ws.Range("A44").Style = "Normal"
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = '9/22/2025'
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").Value = 'AZURDUY JUANA 2449'
$ws.Range("D44").Value = 13
$ws.Range("E44").Value = 'ICD30952422'
$ws.Range("F44").Value = 'Optical Power'
$ws.Range("G44").Value = 'Pendiente'
$ws.Range("H44").Value = 'Cable en panza y cortado'
$ws.Range("I44").Value = 1
$ws.Range("J44").Value = '{"direccionesNormalizadas": [{"altura": 2449, "cod_calle": 1151, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.467279", "y": "-34.551117"}, "direccion": "AZURDUY JUANA 2449, CABA", "nombre_calle": "AZURDUY JUANA", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K44").Value = -58.467279
$ws.Range("L44").Value = -34.551117
$ws.Range("M44").Value = 'Saavedra'
$ws.Range("N44").Value = 'Capital Norte'

# Row 45
$ws.Range("A45").NumberFormat = "@"
$ws.Range("A45").Value = '3069'
$ws.Range("A45").Style = "Normal"
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = '9/22/2025'
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").Value = 'PJE. CAPITAN SAMUEL SPIRO 5996'
$ws.Range("D45").Value = 9
$ws.Range("E45").Value = 'ICD30952704'
$ws.Range("F45").Value = 'Optical Power'
$ws.Range("G45").Value = 'Pendiente'
$ws.Range("H45").Value = 'Cable colgando'
$ws.Range("I45").Value = 1
$ws.Range("J45").Value = '{"direccionesNormalizadas": [{"altura": 5996, "cod_calle": 20114, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.504798", "y": "-34.649012"}, "direccion": "SPIRO, SAMUEL, CAPITAN 5996, CABA", "nombre_calle": "SPIRO, SAMUEL, CAPITAN", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K45").Value = -58.504798
$ws.Range("L45").Value = -34.649012
$ws.Range("M45").Value = 'Devoto'
$ws.Range("N45").Value = 'Capital Norte'

# Row 46
$ws.Range("A46").NumberFormat = "@"
$ws.Range("A46").Value = '3136'
$ws.Range("A46").Style = "Normal"
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = '9/22/2025'
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").Value = 'GAONA AV. 4565'
$ws.Range("D46").Value = 10
$ws.Range("E46").Value = 'ICD30952712'
$ws.Range("F46").Value = 'Optical Power'
$ws.Range("G46").Value = 'Pendiente'
$ws.Range("H46").Value = 'Tendido a Baja altura'
$ws.Range("I46").Value = 1
$ws.Range("J46").Value = '{"direccionesNormalizadas": [{"altura": 4565, "cod_calle": 7025, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.487324", "y": "-34.626415"}, "direccion": "GAONA AV. 4565, CABA", "nombre_calle": "GAONA AV.", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K46").Value = -58.487324
$ws.Range("L46").Value = -34.626415
$ws.Range("M46").Value = 'Devoto'
$ws.Range("N46").Value = 'Capital Norte'

# Row 47
$ws.Range("A47").NumberFormat = "@"
$ws.Range("A47").Value = '3145'
$ws.Range("A47").Style = "Normal"
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = '9/22/2025'
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").Value = 'INCLAN 4329'
$ws.Range("D47").Value = 5
$ws.Range("E47").Value = 'Pendiente ADM'
$ws.Range("F47").Value = 'Optical Power'
$ws.Range("G47").Value = 'Pendiente'
$ws.Range("H47").Value = 'Tendido a baja altura'
$ws.Range("I47").Value = 1
$ws.Range("J47").Value = '{"direccionesNormalizadas": [{"altura": 4329, "cod_calle": 9009, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.424297", "y": "-34.633565"}, "direccion": "INCLAN 4329, CABA", "nombre_calle": "INCLAN", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K47").Value = -58.424297
$ws.Range("L47").Value = -34.633565
$ws.Range("M47").Value = 'Boedo'
$ws.Range("N47").Value = 'Capital Sur'

# Row 48
$ws.Range("A48").NumberFormat = "@"
$ws.Range("A48").Value = '3221'
$ws.Range("A48").Style = "Normal"
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = '9/22/2025'
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = 'GONZALEZ, ELPIDIO 4649'
$ws.Range("D48").Value = 10
$ws.Range("E48").Value = 'ICD30952949'
$ws.Range("F48").Value = 'Optical Power'
$ws.Range("G48").Value = 'Pendiente'
$ws.Range("H48").Value = 'Tendido aereo a baja altura'
$ws.Range("I48").Value = 1
$ws.Range("J48").Value = '{"direccionesNormalizadas": [{"altura": 4649, "cod_calle": 7065, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.501345", "y": "-34.621591"}, "direccion": "GONZALEZ, ELPIDIO 4649, CABA", "nombre_calle": "GONZALEZ, ELPIDIO", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K48").Value = -58.501345
$ws.Range("L48").Value = -34.621591
$ws.Range("M48").Value = 'Devoto'
$ws.Range("N48").Value = 'Capital Norte'

# Row 49
$ws.Range("A49").NumberFormat = "@"
$ws.Range("A49").Value = '3213'
$ws.Range("A49").Style = "Normal"
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = '9/22/2025'
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = 'ALVAREZ JONTE AV. 3705'
$ws.Range("D49").Value = 11
$ws.Range("E49").Value = 'ICD30954655'
$ws.Range("F49").Value = 'Optical Power'
$ws.Range("G49").Value = 'Pendiente'
$ws.Range("H49").Value = 'Cable en panza'
$ws.Range("I49").Value = 1
$ws.Range("J49").Value = '{"direccionesNormalizadas": [{"altura": 3705, "cod_calle": 1056, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.493649", "y": "-34.613319"}, "direccion": "ALVAREZ JONTE AV. 3705, CABA", "nombre_calle": "ALVAREZ JONTE AV.", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K49").Value = -58.493649
$ws.Range("L49").Value = -34.613319
$ws.Range("M49").Value = 'Devoto'
$ws.Range("N49").Value = 'Capital Norte'

# Row 50
$ws.Range("A50").NumberFormat = "@"
$ws.Range("A50").Value = '4249'
$ws.Range("A50").Style = "Normal"
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = '9/22/2025'
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").Value = 'RONDEAU 1501'
$ws.Range("D50").Value = 1
$ws.Range("E50").Value = 'Pendiente ADM'
$ws.Range("F50").Value = 'Optical Power'
$ws.Range("G50").Value = 'Pendiente'
$ws.Range("H50").Value = 'Tendido a baja altura'
$ws.Range("I50").Value = 1
$ws.Range("J50").Value = '{"direccionesNormalizadas": [{"altura": 1501, "cod_calle": 19082, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.386725", "y": "-34.632368"}, "direccion": "RONDEAU 1501, CABA", "nombre_calle": "RONDEAU", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K50").Value = -58.386725
$ws.Range("L50").Value = -34.632368
$ws.Range("M50").Value = 'San Telmo'
$ws.Range("N50").Value = 'Capital Sur'

# Row 51
$ws.Range("A51").NumberFormat = "@"
$ws.Range("A51").Value = '2697'
$ws.Range("A51").Style = "Normal"
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = '9/24/2025'
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Value = 'Mendoza 2299'
$ws.Range("D51").Value = 13
$ws.Range("E51").Value = 'Pendiente ADM'
$ws.Range("F51").Value = 'Optical Power'
$ws.Range("G51").Value = 'Pendiente'
$ws.Range("H51").Value = 'Tendido a baja altura y sin tensar'
$ws.Range("I51").Value = 1
$ws.Range("J51").Value = '{"direccionesNormalizadas": [{"altura": 2299, "cod_calle": 13071, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.456424", "y": "-34.560403"}, "direccion": "MENDOZA 2299, CABA", "nombre_calle": "MENDOZA", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}, {"altura": 2299, "cod_calle": 13072, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.358800", "y": "-34.641733"}, "direccion": "DON PEDRO DE MENDOZA AV. 2299, CABA", "nombre_calle": "DON PEDRO DE MENDOZA AV.", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K51").Value = -58.456424
$ws.Range("L51").Value = -34.560403
$ws.Range("M51").Value = 'Saavedra'
$ws.Range("N51").Value = 'Capital Norte'

# Row 52
$ws.Range("A52").NumberFormat = "@"
$ws.Range("A52").Value = '2725'
$ws.Range("A52").Style = "Normal"
$ws.Range("B52").NumberFormat = "@"
$ws.Range("B52").Value = '9/24/2025'
$ws.Range("B52").Style = "Normal"
$ws.Range("C52").Value = 'QUEVEDO 3352'
$ws.Range("D52").Value = 11
$ws.Range("E52").Value = 'Pendiente ADM'
$ws.Range("F52").Value = 'Optical Power'
$ws.Range("G52").Value = 'Pendiente'
$ws.Range("H52").Value = 'Cable en panza'
$ws.Range("I52").Value = 1
$ws.Range("J52").Value = '{"direccionesNormalizadas": [{"altura": 3352, "cod_calle": 18006, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.518104", "y": "-34.611915"}, "direccion": "QUEVEDO 3352, CABA", "nombre_calle": "QUEVEDO", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K52").Value = -58.518104
$ws.Range("L52").Value = -34.611915
$ws.Range("M52").Value = 'Devoto'
$ws.Range("N52").Value = 'Capital Norte'

# Row 53
$ws.Range("A53").NumberFormat = "@"
$ws.Range("A53").Value = '2729'
$ws.Range("A53").Style = "Normal"
$ws.Range("B53").NumberFormat = "@"
$ws.Range("B53").Value = '9/24/2025'
$ws.Range("B53").Style = "Normal"
$ws.Range("C53").Value = 'SANABRIA 4785'
$ws.Range("D53").Value = 11
$ws.Range("E53").Value = 'Pendiente ADM'
$ws.Range("F53").Value = 'Optical Power'
$ws.Range("G53").Value = 'Pendiente'
$ws.Range("H53").Value = 'Cables cortados'
$ws.Range("I53").Value = 1
$ws.Range("J53").Value = '{"direccionesNormalizadas": [{"altura": 4785, "cod_calle": 20021, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.520172", "y": "-34.596414"}, "direccion": "SANABRIA 4785, CABA", "nombre_calle": "SANABRIA", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K53").Value = -58.520172
$ws.Range("L53").Value = -34.596414
$ws.Range("M53").Value = 'Paternal'
$ws.Range("N53").Value = 'Capital Norte'

# Row 54
$ws.Range("A54").NumberFormat = "@"
$ws.Range("A54").Value = '2921'
$ws.Range("A54").Style = "Normal"
$ws.Range("B54").NumberFormat = "@"
$ws.Range("B54").Value = '9/24/2025'
$ws.Range("B54").Style = "Normal"
$ws.Range("C54").Value = 'SARMIENTO 4270'
$ws.Range("D54").Value = 5
$ws.Range("E54").Value = 'Pendiente ADM'
$ws.Range("F54").Value = 'Optical Power'
$ws.Range("G54").Value = 'Pendiente'
$ws.Range("H54").Value = 'Cable en panza'
$ws.Range("I54").Value = 1
$ws.Range("J54").Value = '{"direccionesNormalizadas": [{"altura": 4270, "cod_calle": 20074, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.425368", "y": "-34.604407"}, "direccion": "SARMIENTO 4270, CABA", "nombre_calle": "SARMIENTO", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}, {"altura": 4270, "cod_calle": 20075, "cod_calle_cruce": null, "cod_partido": "caba", "coordenadas": {"srid": 4326, "x": "-58.410672", "y": "-34.570296"}, "direccion": "SARMIENTO AV. 4270, CABA", "nombre_calle": "SARMIENTO AV.", "nombre_calle_cruce": "", "nombre_localidad": "CABA", "nombre_partido": "CABA", "tipo": "calle_altura"}]}'
$ws.Range("K54").Value = -58.425368
$ws.Range("L54").Value = -34.604407
$ws.Range("M54").Value = 'Almagro'
$ws.Range("N54").Value = 'Capital Sur'

